$d = $word.ActiveDocument

$d.Content.Find.Execute("21. - -Astudy was conducted fo identify the factors that affect the evaporation rate", $true, $false, $false, $false, $false, $true, 1, $false, "21. - Astudy was conducted fo identify the factors that affect the evaporation rate", 2)

$d.Content.Find.Execute("-: of:some liquids. The table below shows the results.of the investigation.", $true, $false, $false, $false, $false, $true, 1, $false, "-. of some liquids. The table below shows the results of the investigation.", 2)

$d.Content.Find.Execute("ema ea", $true, $false, $false, $false, $false, $true, 1, $false, "emt ra a", 2)

$d.Content.Find.Execute("Wate [8 [eo [ease [92 [et |", $true, $false, $false, $false, $false, $true, 1, $false, "Wate [8 [eo [ea [se [92 [er |", 2)

$d.Content.Find.Execute("(1) The evaporation rate for water is less than that for alcohol.", $true, $false, $false, $false, $false, $true, 1, $false, "(1) |The evaporation rate tor waier is tess inan tnat for aiconol.", 2)

$d.Content.Find.Execute("(2) The larger the amounts of water the higher the evaporation rate.", $true, $false, $false, $false, $false, $true, 1, $false, "(2} The larger the amounts of water the higher the evaporation rate.", 2)
